$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update standard-error cells (shown in parentheses) under the
# theta / lambda / proportion-drinking columns for the new bootstrapping.
$ws.Range("C3").Value = "(0.7)"
$ws.Range("D3").Value = "(0.3)"

$ws.Range("C5").Value = "(0.28)"
$ws.Range("D5").Value = "(0.11)"

$ws.Range("C7").Value = "(0.46)"
$ws.Range("D7").Value = "(0.44)"

$ws.Range("C9").Value = "(0.17)"
$ws.Range("D9").Value = "(0.26)"

$ws.Range("C11").Value = "(0.39)"
$ws.Range("D11").Value = "(0.54)"

$ws.Range("C13").Value = "(0.18)"
$ws.Range("D13").Value = "(0.48)"

$ws.Range("C15").Value = "(1.91)"
$ws.Range("D15").Value = "(1.16)"
